$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 3 (the old row 3 - speedometer_text - moves to row 4)
$ws.Rows.Item(3).Insert()

# The inserted row inherits formatting (and stray styled-but-empty cells) from
# the row above across its full column span; clear those out so only A:D are populated
$ws.Range("E3:H3").Clear()

# Fill in the new row 3 with the recurring-challenge text-block template entry
# (values assigned in this order so new shared strings land in the same
# table order as the target workbook: description, variable name, template)
$ws.Range("A3").Value = "Summary_Report"
$ws.Range("C3").Value = "One feature of the auto-generated template is that it pulls the most common recurring challenges across the agency. This text block summarizes how many times a recurring challenge has occurred for a given agency."
$ws.Range("B3").Value = "recurring_challenge_text"
$ws.Range("D3").Value = "**{challenge}** has been reported as challenge for the **{goal}** team in each of the last **{challenge count} quarters**."

# Match formatting (wrap text) used by the rest of the data rows
$ws.Range("A3:D3").WrapText = $true

# Row heights: new row taller (wraps more text), old row (now row 4) keeps its height
$ws.Rows.Item(3).RowHeight = 187.2
$ws.Rows.Item(4).RowHeight = 86.4

# Update the active selection to D3, matching the saved view state
$ws.Range("D3").Select()
